# Add a new IHU requirement (REQ-076: "The IHU shall operate as a RTOS.")
# The new requirement's Description/Reasoning/Priority Level are inserted into
# row 59 (just after REQ-055), pushing the existing Description/Reasoning/
# Priority Level values of rows 59-78 down by one row. The Reference Number
# (col A) and Status (col E) columns are left as a simple, already-contiguous
# sequence, so a brand new row 79 is appended at the bottom carrying the new
# REQ-076 reference number together with the Description/Reasoning/Priority
# that used to belong to the last row (REQ-075).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create row 79 by copying the formatting (fonts, wrap text, etc.) of the
# current last data row (78) and inserting it below, shifting nothing else.
$ws.Rows.Item(78).Copy()
$ws.Rows.Item(79).Insert(-4121)
$excel.CutCopyMode = 0

# Shift the Description (B), Reasoning (C) and Priority Level (D) values for
# rows 59-78 down to rows 60-79, working bottom-up so we never clobber a
# value before it has been copied.
for ($r = 78; $r -ge 59; $r--) {
    $destRow = $r + 1
    $ws.Range("B$destRow").Value = $ws.Range("B$r").Value()
    $ws.Range("C$destRow").Value = $ws.Range("C$r").Value()
    $ws.Range("D$destRow").Value = $ws.Range("D$r").Value()
}

# Fill in the new requirement's Description/Reasoning/Priority on row 59.
$ws.Range("B59").Value = "The IHU shall operate as a RTOS."
$ws.Range("C59").Value = "Need to process mission critical events as they happen."
$ws.Range("D59").Value = "High"

# The Reference Number and Status columns simply gain one new contiguous
# entry at the very end of the table.
$ws.Range("A79").Value = "REQ-076"
$ws.Range("E79").Value = "Not Done"
